# Adds a "TabName" column (A) identifying this row as the "CasesTab",
# and rewrites the Neo4j queries to also pull gender/ethnicity (TC for CTDC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:D columns to B:E, making room for the new TabName column.
$ws.Columns("A").Insert()

# New column A: tab identifier.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Columns("A").ColumnWidth = 8

# Updated queries (now also returning Gender / Ethnicity) in the shifted B/C columns.
$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.ethnicity IN ['UNKNOWN']`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"
$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE WHERE c.ethnicity IN ['UNKNOWN']`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# The row grew taller to fit the longer wrapped query text.
$ws.Rows(2).RowHeight = 174

$ws.Range("B2").Select() | Out-Null
